$wb = $excel.ActiveWorkbook

# --- Update the SYP description text on Sheet3 (B30): SYP -> Y wording ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B30").Value = "J = junior, '' = standard, S = senior, Y = student/young professional"

# Move the cursor/selection on Sheet3 before switching away (matches diff: E16 -> C12)
$ws3.Activate()
$ws3.Range("C12").Select() | Out-Null

# --- Add the new Sheet4 with club/database mapping data ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"

$ws4.Range("A1").Value = "Club Name"
$ws4.Range("B1").Value = "Club ID"
$ws4.Range("C1").Value = "Database"

$ws4.Range("A2").Value = "Highpoint"
$ws4.Range("B2").Value = 201
$ws4.Range("C2").Value = "NMSW"

$ws4.Range("A3").Value = "Midtown"
$ws4.Range("B3").Value = 202
$ws4.Range("C3").Value = "NMSW"

$ws4.Range("A4").Value = "Downtown"
$ws4.Range("B4").Value = 203
$ws4.Range("C4").Value = "NMSW"

$ws4.Range("A5").Value = "Del Norte"
$ws4.Range("B5").Value = 204
$ws4.Range("C5").Value = "NMSW"

$ws4.Range("A6").Value = "Riverpoint"
$ws4.Range("B6").Value = 205
$ws4.Range("C6").Value = "NMSW"

$ws4.Range("A7").Value = "DTC"
$ws4.Range("B7").Value = 252
$ws4.Range("C7").Value = "Denver"

$ws4.Range("A8").Value = "Tabor Center"
$ws4.Range("B8").Value = 254
$ws4.Range("C8").Value = "Denver"

$ws4.Range("A9").Value = "Flatirons"
$ws4.Range("B9").Value = 257
$ws4.Range("C9").Value = "Denver"

$ws4.Range("A10").Value = "Monaco"
$ws4.Range("B10").Value = 292
$ws4.Range("C10").Value = "Denver"

$ws4.Range("A11").Value = "MAC"
$ws4.Range("B11").Value = 375
$ws4.Range("C11").Value = "MAC"

# Column A best-fit width (matches diff col def width="12.42578125" bestFit="1")
$ws4.Columns.Item(1).AutoFit() | Out-Null

$ws4.Range("A1:C11").Select() | Out-Null
$ws4.Activate()
